$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for Wins, Losses, Ties in AD1:AF1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$hdrRange = $ws.Range("AD1:AF1")
$hdrRange.Font.Bold = $true
$hdrRange.HorizontalAlignment = -4108
$hdrRange.VerticalAlignment = -4160
$hdrRange.Borders.LineStyle = 1

# Fill in team record (Wins/Losses/Ties) for every player data row (2-49)
for ($r = 2; $r -le 49; $r++) {
    $ws.Cells.Item($r, 30).Value = 56
    $ws.Cells.Item($r, 31).Value = 106
    $ws.Cells.Item($r, 32).Value = 0
}

Write-Output "done"
